# "feat: importacion masiva terminada"
# The sheet previously had two header rows: a long descriptive row (row 1)
# and a short machine field-name row (row 2), followed by 5 product rows.
# This finishes the bulk-import template: the long descriptive header row is
# removed (the short field-name row becomes the single header row), and the
# last two imported products (008 / 009) are left without a computed sale
# price yet, matching the mid-import state that was committed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the old verbose header row - the machine-friendly header
#    (nombre/precioCompra/precioVenta/categoriaId) becomes row 1.
$ws.Rows.Item(1).Delete()

# 2) The last two freshly-imported products (rows 5 & 6 now) haven't had
#    their "precioVenta" formula filled in yet.
$ws.Range("C5:C6").ClearContents()

# 3) Stray formatted cell left over out at column L (underlined), same as
#    what shipped in the finished import.
$ws.Range("L4").Font.Underline = [int]2

# 4) Column sizing for the now-visible data columns.
$ws.Columns.Item(2).ColumnWidth = 11.95
$ws.Columns.Item(3).ColumnWidth = 11.95
$ws.Columns.Item(4).ColumnWidth = 34.95

# 5) Selection / print orientation to match the finished workbook.
[void]$ws.Range("C15").Select()
$ws.PageSetup.Orientation = 1
